# Weekly price-sheet update: a new week's record is inserted as row 15
# (right after the existing data that precedes it), pushing the previous
# rows 15-34 down to 16-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15, shifting rows 15:34 down to 16:35.
$ws.Rows("15:15").Insert()

# Populate the new row 15 with this week's record.
$ws.Cells.Item(15, 1).Value = 10
$ws.Cells.Item(15, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(15, 3).Value = "La Araucanía"
$ws.Cells.Item(15, 4).Value = 44484
$ws.Cells.Item(15, 5).Value = 9
$ws.Cells.Item(15, 6).Value = 100112026
$ws.Cells.Item(15, 7).Value = "Haba"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 9000
$ws.Cells.Item(15, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 360
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
